$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Thu Feb 16 12:11:44 EST 2023"
$ws.Range("B3").Value = "Thu Feb 16 12:11:54 EST 2023"
$ws.Range("B4").Value = "Thu Feb 16 12:12:04 EST 2023"
